# Applies the benchmark-stats update to the single-column results table.
$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Simple value replacements in the first 6 rows (indices stable) ---
$t.Cell(1, 1).Range.Text = "0M"        # was 99.96
$t.Cell(2, 1).Range.Text = "0M"        # was 0.28
$t.Cell(3, 1).Range.Text = "0M"        # was 752
$t.Cell(4, 1).Range.Text = "1856"      # was 826
$t.Cell(5, 1).Range.Text = "0.00001"   # was 0.00003
# row 6 (0.00239) is left unchanged

# --- Remove rows 7,8,9 (0.00010 / 0.00007 / 0.00011) entirely ---
# Delete from the highest index down so the remaining indices don't shift.
$t.Rows(9).Delete()
$t.Rows(8).Delete()
$t.Rows(7).Delete()

# After the deletions, the former row 10 (0.00013) is now row 7 (unchanged),
# former row 11 (0.00017) is now row 8, former row 12 (0.09565) is now row 9.
$t.Cell(8, 1).Range.Text = "0.00006"   # was 0.00017
$t.Cell(9, 1).Range.Text = "0.00020"   # was 0.09565

# --- Insert three new rows right after the row that now holds 0.00020,
#     and before the row holding 100.0 ---
# Each Rows.Add(beforeRow) call inserts immediately above "beforeRow", so to
# end up with 0.00025, 0.00030, 0.28238 (in that reading order) the rows must
# be added in the reverse order.
$beforeRow = $t.Rows(10)   # the "100.0" row
$newRow3 = $t.Rows.Add($beforeRow)
$newRow3.Cells(1).Range.Text = "0.28238"
$newRow2 = $t.Rows.Add($beforeRow)
$newRow2.Cells(1).Range.Text = "0.00030"
$newRow1 = $t.Rows.Add($beforeRow)
$newRow1.Cells(1).Range.Text = "0.00025"

# --- Collapse the three multi-tab summary rows (now the last three rows of
#     the table) down to their single headline value ---
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.96"   # was 515 <tab> ... <tab> 100.0
$t.Cell($rowCount - 1, 1).Range.Text = "0.28"    # was 204 <tab> ... <tab> 100.0
$t.Cell($rowCount, 1).Range.Text = "752"         # was 311 <tab> ... <tab> 100.0
